$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 443
$ws.Cells.Item(5, 6).Value = 1826
$ws.Cells.Item(7, 6).Value = 2899
$ws.Cells.Item(8, 6).Value = 2417
$ws.Cells.Item(9, 6).Value = 758
$ws.Cells.Item(10, 6).Value = 7309
$ws.Cells.Item(11, 6).Value = 186
$ws.Cells.Item(13, 6).Value = 194
$ws.Cells.Item(14, 6).Value = 1669
$ws.Cells.Item(15, 6).Value = 1423
$ws.Cells.Item(16, 6).Value = 1262
$ws.Cells.Item(17, 6).Value = 131
$ws.Cells.Item(18, 6).Value = 131
$ws.Cells.Item(19, 6).Value = 3348
$ws.Cells.Item(20, 6).Value = 5782
$ws.Cells.Item(21, 6).Value = 5782
$ws.Cells.Item(22, 6).Value = 575
$ws.Cells.Item(23, 6).Value = 938
$ws.Cells.Item(24, 6).Value = 1181
$ws.Cells.Item(25, 6).Value = 337
$ws.Cells.Item(26, 6).Value = 5751
$ws.Cells.Item(27, 6).Value = 323
$ws.Cells.Item(28, 6).Value = 49
$ws.Cells.Item(29, 6).Value = 3962
$ws.Cells.Item(30, 6).Value = 209
$ws.Cells.Item(31, 6).Value = 661
$ws.Cells.Item(32, 6).Value = 1820
$ws.Cells.Item(33, 6).Value = 1117
$ws.Cells.Item(34, 6).Value = 244
$ws.Cells.Item(35, 6).Value = 14
$ws.Cells.Item(36, 6).Value = 147
$ws.Cells.Item(37, 6).Value = 85
$ws.Cells.Item(38, 6).Value = 301
$ws.Cells.Item(39, 6).Value = 1100
$ws.Cells.Item(41, 6).Value = 1810
$ws.Cells.Item(42, 6).Value = 73
$ws.Cells.Item(43, 6).Value = 343
$ws.Cells.Item(44, 6).Value = 126
$ws.Cells.Item(45, 6).Value = 1002
$ws.Cells.Item(46, 6).Value = 541
$ws.Cells.Item(50, 6).Value = 132
$ws.Cells.Item(51, 6).Value = 8

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 160
$ws.Cells.Item(7, 6).Value = 113
$ws.Cells.Item(12, 6).Value = 11
$ws.Cells.Item(13, 6).Value = 84
$ws.Cells.Item(14, 6).Value = 640
$ws.Cells.Item(15, 6).Value = 317
$ws.Cells.Item(22, 6).Value = 331
$ws.Cells.Item(27, 6).Value = 64
$ws.Cells.Item(35, 6).Value = 10

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 3321
$ws.Cells.Item(5, 6).Value = 429
$ws.Cells.Item(7, 6).Value = 1518
$ws.Cells.Item(8, 6).Value = 774
$ws.Cells.Item(9, 6).Value = 445
$ws.Cells.Item(10, 6).Value = 2989
$ws.Cells.Item(11, 6).Value = 384
$ws.Cells.Item(12, 6).Value = 745
$ws.Cells.Item(13, 6).Value = 927
$ws.Cells.Item(14, 6).Value = 944
$ws.Cells.Item(15, 6).Value = 1420

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 429
$ws.Cells.Item(3, 6).Value = 774
$ws.Cells.Item(5, 6).Value = 443
$ws.Cells.Item(6, 6).Value = 445
$ws.Cells.Item(7, 6).Value = 2989
$ws.Cells.Item(9, 6).Value = 2899
$ws.Cells.Item(10, 6).Value = 2417
$ws.Cells.Item(11, 6).Value = 758
$ws.Cells.Item(12, 6).Value = 7309
$ws.Cells.Item(13, 6).Value = 186
$ws.Cells.Item(14, 6).Value = 745
$ws.Cells.Item(16, 6).Value = 194
$ws.Cells.Item(17, 6).Value = 1423
$ws.Cells.Item(18, 6).Value = 944
$ws.Cells.Item(19, 6).Value = 640
$ws.Cells.Item(20, 6).Value = 131
$ws.Cells.Item(21, 6).Value = 3348
$ws.Cells.Item(22, 6).Value = 317
$ws.Cells.Item(23, 6).Value = 5782
$ws.Cells.Item(24, 6).Value = 575
$ws.Cells.Item(25, 6).Value = 938
$ws.Cells.Item(26, 6).Value = 1181
$ws.Cells.Item(27, 6).Value = 337
$ws.Cells.Item(28, 6).Value = 5751
$ws.Cells.Item(29, 6).Value = 323
$ws.Cells.Item(30, 6).Value = 3962
$ws.Cells.Item(31, 6).Value = 661
$ws.Cells.Item(32, 6).Value = 331
$ws.Cells.Item(33, 6).Value = 1821
$ws.Cells.Item(34, 6).Value = 1117
$ws.Cells.Item(36, 6).Value = 147
$ws.Cells.Item(37, 6).Value = 85
$ws.Cells.Item(38, 6).Value = 301
$ws.Cells.Item(39, 6).Value = 1100
$ws.Cells.Item(40, 6).Value = 1810
$ws.Cells.Item(41, 6).Value = 73
$ws.Cells.Item(42, 6).Value = 343
$ws.Cells.Item(43, 6).Value = 126
$ws.Cells.Item(44, 6).Value = 1002
$ws.Cells.Item(46, 6).Value = 541
$ws.Cells.Item(50, 6).Value = 132
$ws.Cells.Item(51, 6).Value = 8
